# Update the "dSF" (column F) values for the rows whose underlying data was
# repulled. This reflects the commit "repull data, push all data, mean
# calculation" which corrected several dSF values while leaving the other
# columns (including dS0 / column E) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -4
    12 = -2
    19 = -3
    22 = 0
    24 = -1
    25 = -4
    33 = 4
    37 = -6
    44 = -1
    47 = -2
    50 = -3
    51 = -1
    53 = -3
    56 = -2
    59 = -6
    60 = -3
    62 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
